# Verified and checked matchmaking algorithm
#
# Rebuilds the "Sheet2" worksheet (the order-book simulation output) with
# the fuller A:L layout: pre-trade orderbook (rows 2-11), trades
# (rows 12-25) and post-trade orderbook (rows 28-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Start from a clean sheet - the previous layout (B3:G11 + C16:E27) is
# being fully replaced.
$ws.Cells.Clear()

# ------------------------------------------------------------------
# Section header: pre-trade orderbook
# ------------------------------------------------------------------
$ws.Range("A2").Value = "pre-trader orderbook"

# Pre-trade orderbook rows (Asks = "A(uuid", Bids = "B(uuid")
$preTrade = @(
    @{ row=3;  label="A(uuid"; orderid=1; vol=61; price=70; sign=-1 },
    @{ row=4;  label="A(uuid"; orderid=2; vol=50; price=50; sign=-1 },
    @{ row=5;  label="A(uuid"; orderid=3; vol=40; price=35; sign=-1 },
    @{ row=6;  label="A(uuid"; orderid=4; vol=45; price=35; sign=-1 },
    @{ row=7;  label="B(uuid"; orderid=0; vol=20; price=29; sign=1  },
    @{ row=8;  label="B(uuid"; orderid=0; vol=47; price=31; sign=1  },
    @{ row=9;  label="B(uuid"; orderid=0; vol=43; price=32; sign=1  },
    @{ row=10; label="B(uuid"; orderid=0; vol=33; price=31; sign=1  }
)

foreach ($o in $preTrade) {
    $r = $o.row
    $ws.Range("B$r").Value = $o.label
    $ws.Range("C$r").Value = $o.orderid
    $ws.Range("D$r").Value = " orderid"
    if ($o.sign -lt 0) {
        $ws.Range("E$r").Formula = "=-G$r"
    } else {
        $ws.Range("E$r").Formula = "=G$r"
    }
    $ws.Range("F$r").Value = " volume"
    $ws.Range("G$r").Value = $o.vol
    $ws.Range("H$r").Value = " price"
    $ws.Range("I$r").Value = $o.price
    $ws.Range("J$r").Value = " timeStamp"
    if ($r -eq 10) {
        $ws.Range("K$r").Value = " 1)]"
    } else {
        $ws.Range("K$r").Value = " 1)"
        $ws.Range("L$r").Value = " "
    }
}

$ws.Range("E11").Formula = "=SUM(E3:E10)"

# ------------------------------------------------------------------
# Section header: trades
# ------------------------------------------------------------------
$ws.Range("A12").Value = "trades"

$trades = @(
    @{ row=13; orderid=1; typ=" OrderType.ASK"; vol=20; priceLabel=" 49.5)" },
    @{ row=14; orderid=1; typ=" OrderType.ASK"; vol=41; priceLabel=" 50.5)" },
    @{ row=15; orderid=2; typ=" OrderType.ASK"; vol=39; priceLabel=" 40.5)" },
    @{ row=16; orderid=2; typ=" OrderType.ASK"; vol=11; priceLabel=" 41.0)" },
    @{ row=17; orderid=4; typ=" OrderType.ASK"; vol=17; priceLabel=" 33.5)" },
    @{ row=18; orderid=3; typ=" OrderType.ASK"; vol=15; priceLabel=" 33.5)" },
    @{ row=19; orderid=0; typ=" OrderType.BID"; vol=20; priceLabel=" 49.5)" },
    @{ row=20; orderid=0; typ=" OrderType.BID"; vol=17; priceLabel=" 50.5)" },
    @{ row=21; orderid=0; typ=" OrderType.BID"; vol=24; priceLabel=" 50.5)" },
    @{ row=22; orderid=0; typ=" OrderType.BID"; vol=16; priceLabel=" 40.5)" },
    @{ row=23; orderid=0; typ=" OrderType.BID"; vol=23; priceLabel=" 40.5)" },
    @{ row=24; orderid=0; typ=" OrderType.BID"; vol=11; priceLabel=" 41.0)" },
    @{ row=25; orderid=0; typ=" OrderType.BID"; vol=32; priceLabel=" 33.5)]" }
)

foreach ($t in $trades) {
    $r = $t.row
    $ws.Range("B$r").Value = "T(uuid"
    $ws.Range("C$r").Value = $t.orderid
    $ws.Range("D$r").Value = " order id"
    if ($t.typ -eq " OrderType.ASK") {
        $ws.Range("E$r").Formula = "=-I$r"
    } else {
        $ws.Range("E$r").Formula = "=I$r"
    }
    $ws.Range("F$r").Value = " Type"
    $ws.Range("G$r").Value = $t.typ
    $ws.Range("H$r").Value = " Volume"
    $ws.Range("I$r").Value = $t.vol
    $ws.Range("J$r").Value = " Price"
    $ws.Range("K$r").Value = $t.priceLabel
    if ($r -ne 25) {
        $ws.Range("L$r").Value = " "
    }
}

# ------------------------------------------------------------------
# Section header: post-trade orderbook
# ------------------------------------------------------------------
$ws.Range("A28").Value = "post-trade orderbook"
$ws.Range("B29").Value = "Orderbook"

$postTrade = @(
    @{ row=30; orderid=4; vol=28; price=35 },
    @{ row=31; orderid=3; vol=25; price=35 }
)

foreach ($o in $postTrade) {
    $r = $o.row
    $ws.Range("B$r").Value = "A(uuid"
    $ws.Range("C$r").Value = $o.orderid
    $ws.Range("D$r").Value = " orderid"
    $ws.Range("E$r").Value = 0
    $ws.Range("F$r").Value = " volume"
    $ws.Range("G$r").Value = $o.vol
    $ws.Range("H$r").Value = " price"
    $ws.Range("I$r").Value = $o.price
    $ws.Range("J$r").Value = " timeStamp"
    if ($r -eq 31) {
        $ws.Range("K$r").Value = " 1)]"
    } else {
        $ws.Range("K$r").Value = " 1)"
        $ws.Range("L$r").Value = " "
    }
}

# ------------------------------------------------------------------
# Column widths: drop the old single wide column F, replace with a
# narrower F and a new G.
# ------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 10.166666666666666
$ws.Columns.Item(7).ColumnWidth = 28.307291666666668

# ------------------------------------------------------------------
# Re-apply the sort on the trades block (by column G, i.e. order Type)
# so the worksheet keeps a sortState matching the on-screen order.
# ------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("G13:G25")) | Out-Null
$ws.Sort.SetRange($ws.Range("A13:L26"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# ------------------------------------------------------------------
# Selection / view state
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("G29:G33").Select()
